$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The placeholder tick-mark in N30 becomes "x"
$ws.Range("N30").Value = "x"

# SQL comment lines added to column I, two blocks of commented-out SQL
# (the original author's scratch work for the "assignment").
$lines = @{
  27 = '-- SELECT i.industry_name'
  28 = '-- FROM industry as i '
  29 = '-- INNER JOIN series as s '
  30 = '-- ON i.industry_code = s.industry_code'
  31 = '-- SELECT MAX(j.value) as "Highest Average Weekly Hours Industry", MIN(j.value) as "Lowest Average Weekly Hours Industry", i.industry_name'
  32 = '-- FROM industry as i '
  33 = '-- JOIN series as s '
  34 = '-- ON i.industry_code = s.industry_code'
  35 = '-- JOIN january_2017 as j '
  36 = '-- ON j.series_id = s.series_id'
  42 = '-- SELECT MAX(j.value) as "Highest Average Weekly Hours Industry", MIN(j.value) as "Lowest Average Weekly Hours Industry", i.industry_name'
  43 = '-- FROM january_2017 as j'
  44 = '-- INNER JOIN series as s '
  45 = '-- ON j.series_id = s.series_id '
  46 = '-- INNER JOIN datatype as d '
  47 = '-- ON d.data_type_code = s.data_type_code '
  48 = '-- INNER JOIN industry as i '
  49 = '-- ON i.industry_code = s.industry_code'
  50 = '-- WHERE j.series_id IN (SELECT s.series_id from series as s WHERE s.data_type_code IN '
  51 = "-- (SELECT d.data_type_code from datatype as d WHERE d.data_type_text = 'Average weekly hours of production and nonsupervisory employees'))"
  52 = '-- GROUP BY i.industry_name'
  53 = '-- ORDER BY MAX(j.value) ASC'
  54 = '-- WHERE j.series_id IN (SELECT s.series_id from series as s WHERE s.data_type_code IN '
  55 = "-- (SELECT d.data_type_code from datatype as d WHERE d.data_type_text = 'Average weekly hours of production and nonsupervisory employees'))"
  56 = '-- GROUP BY i.industry_name'
  57 = '-- HAVING MAX(j.value), MIN(j.value)'
}

$rows = 27..36 + 42..57

# Build the small green "code comment" font/style once, on the first row,
# then stamp the same format onto every other commented line via a
# format-only paste so they all share one cell style.
$first = $ws.Range("I27")
$first.Value = $lines[27]
$first.Font.Name = "Consolas"
$first.Font.Size = 7
$first.Font.Color = 32768

$first.Copy() | Out-Null
foreach ($r in $rows) {
  if ($r -eq 27) { continue }
  $cell = $ws.Range("I$r")
  $cell.Value = $lines[$r]
  $cell.PasteSpecial(-4122) | Out-Null
}

# Match the selection left behind by the edit
$ws.Range("I42:I57").Select() | Out-Null
